$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 24
    7  = 86
    8  = 456
    9  = 48
    10 = 22
    11 = 568
    12 = 29
    15 = 370
    21 = 95
    22 = 922
    23 = 1395
    24 = 299
    27 = 73
    31 = 214
    34 = 1616
    37 = 157
    40 = 3625
    41 = 424
    42 = 200
    46 = 66
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
